$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

# --- Edit 1 ---------------------------------------------------------------
# "...nach dem Schema Extract, Transform..." -> "...nach dem ETL-Schema, Extract, Transform..."
$d.Content.Find.Execute(
    "Schema Extract, Transform",
    $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "ETL-Schema, Extract, Transform",
    $wdReplaceAll) | Out-Null

# --- Edit 2 ---------------------------------------------------------------
# "In der Extract Phase werden die benötigten Daten für die Problemstellung
#  gesammelt." -> "In der Extract Phase werden Daten aus Quellen im Internet
#  gesammelt."
$d.Content.Find.Execute(
    "In der Extract Phase werden die benötigten Daten für die Problemstellung gesammelt.",
    $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "In der Extract Phase werden Daten aus Quellen im Internet gesammelt.",
    $wdReplaceAll) | Out-Null

# --- Edit 3 ---------------------------------------------------------------
# "...ebscappers von den jeweiligen Webpages heruntergeladen." -> "...ebscappers
#  von den jeweiligen Webpage HTML's extrahiert."
$d.Content.Find.Execute(
    "ebscappers von den jeweiligen Webpages heruntergeladen.",
    $false, $false, $false, $false, $false,
    $true, $wdFindContinue, $false,
    "ebscappers von den jeweiligen Webpage HTML’s extrahiert.",
    $wdReplaceAll) | Out-Null
